$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0.027
$ws.Range("F4").Value = 1.046
$ws.Range("F5").Value = 1.085
$ws.Range("F6").Value = 0.61
$ws.Range("F7").Value = 0.532
$ws.Range("F8").Value = 0.124
$ws.Range("F9").Value = 0.084
$ws.Range("F10").Value = 0.058
$ws.Range("F11").Value = 0.05

# Cells that previously used the tiny 6.5pt font now get the normal
# 10.5pt Times New Roman font used elsewhere in the table, now that
# they hold real values.
$ws.Range("F3,F5,F7,F9,F11").Font.Size = 10.5
$ws.Range("F3,F5,F7,F9,F11").Font.Name = "Times New Roman"

$ws.Application.ActiveWindow.RangeSelection
$ws.Range("F9").Select()
